$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "YES"
$ws.Range("B8").Value = "URL"
$ws.Range("C8").Value = "www.magenta.ca|order"

$ws.Range("C8").Select()
